$d = $word.ActiveDocument

# The "Etapes pour lancer l'application" numbered list originally had these
# three items (numId 2):
#   1. Exécuter le script en tant qu'administrateur *installScript.bat*
#   2. Lancer l'installateur *openimuInstaller.exe*
#   3. Lancer l'exécutable ... applicationOpenimu.exe
#
# The edit drops item 1 entirely, so the list becomes:
#   1. Lancer l'installateur *openimuInstaller.exe*
#   2. Lancer l'exécutable ... applicationOpenimu.exe
#
# It also relocates the hidden "_GoBack" bookmark (which Word stamps at the
# last edit position) from the trailing empty paragraph to the start of the
# new first list item.

# Find the paragraph that still reads "Exécuter le script en tant
# qu'administrateur installScript.bat" and remove it completely, including
# its paragraph mark, so the following item moves up to take its place.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Ex*cuter le script*installScript.bat*") {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}

# Locate the (now first) "Lancer l'installateur openimuInstaller.exe" item.
$installerPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Lancer l*installateur*openimuInstaller.exe*") {
        $installerPara = $p
        break
    }
}

# Move/recreate the hidden "_GoBack" bookmark at the very start of that
# paragraph (adding a bookmark with an existing name relocates it, removing
# the old one automatically).
if ($installerPara -ne $null) {
    $startRange = $d.Range($installerPara.Range.Start, $installerPara.Range.Start)
    $startRange.Bookmarks.Add("_GoBack")
}
